# Apply the commit's changes to the "EZ Water Adjustment" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EZ Water Adjustment")

# --- Core data/formula fix -------------------------------------------------
# The dilution factor in the "Lactic Acid (ml):" helper cell was inverted;
# it should divide, not multiply. Everything downstream (RA, mash/sparge
# CaCO3 figures, Cl:SO4 ratio, and the "Raw Text Format" sheet that mirrors
# these cells) recalculates automatically from this single change.
$ws.Range("E37").Formula = "=2.42/1.335"

# --- Print area --------------------------------------------------------
$ws.PageSetup.PrintArea = "`$A`$1:`$K`$55"

# --- Page scaling: "fit to page" flag + updated scale percentage ----------
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.Zoom = 61
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1

# --- Selection / scroll position -------------------------------------------
# Previously the view was scrolled to A7 with E26 selected; now the view is
# back at the top with the title row (C1:K1) selected.
$ws.Range("C1:K1").Select()
